$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new rule row ("Set Consultation Type") above the
#     existing "Set Consultation Priority" row (old row 25, now row 26). ---
$ws.Rows("25:25").Insert(-4121)

# Re-create the thin-box border formatting used by the rest of the
# rule table (matches style used by B:D in the surrounding rows).
$ws.Range("B25:D25").Borders.LineStyle = 1
$ws.Range("B25:D25").Borders.Weight = 2

$ws.Range("B25").Value = "Set Consultation Type"
$ws.Range("D25").Value = "setConsultationType, 'Consultation'"
$ws.Range("C25").Value = "consultationType == null || consultationType.equals("""")"

# --- Append a new blank row at the bottom of the table (old last row
#     33 becomes 34, and a fresh blank row is appended after it). ---
$ws.Rows("33:33").Copy()
$ws.Rows("34:34").Insert(-4121)

# --- Restore selection/view state recorded for the sheet. ---
$ws.Activate()
$ws.Range("C28").Select()
